# Insert a new row above the existing "Width:" row (row 8) on the first
# worksheet ("Survey 1") and populate it with a new
# "Pseudo-Random Question Width:" label + an (empty) input cell, mirroring
# the layout of the other label/value rows in the survey-table-properties
# block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift row 8 (and everything below it) down by one row.
$ws.Rows.Item(8).Insert()

# New label cell, formatted like the other "value" cells (style copied from
# the adjacent empty input cell B8, which inherited its formatting from the
# old row 8 during the insert).
$ws.Range("A8").Value = "Pseudo-Random Question Width:"
$ws.Range("A8").Style = $ws.Range("B8").Style

# Match the author's recorded selection after making the edit.
$ws.Range("B8").Select()
